$wb = $excel.ActiveWorkbook

# --- Sheet "parameters": update vehicle_capacity and no_stops values ---
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Range("B3").Value = 4
$wsParams.Range("B12").Value = 6
$wsParams.Range("C15").Select()

# --- Sheet "comp_quantity_inst1": update C5, replace rows 7-30 with new 4-row block ---
$wsComp = $wb.Worksheets.Item("comp_quantity_inst1")
$wsComp.Range("C5").Value = 28

# Remove old rows 11-30 entirely, shifting cells up, so dimension shrinks to A1:D10
$wsComp.Range("A11:D30").Delete(-4162)

# Overwrite rows 7-10 with the new data block
$wsComp.Range("A7").Value = "T1"
$wsComp.Range("B7").Value = "T2"
$wsComp.Range("C7").Value = 24
$wsComp.Range("D7").Value = 0

$wsComp.Range("A8").Value = "T2"
$wsComp.Range("B8").Value = "T3"
$wsComp.Range("C8").Value = 322
$wsComp.Range("D8").Value = 0

$wsComp.Range("A9").Value = "T3"
$wsComp.Range("B9").Value = "T4"
$wsComp.Range("C9").Value = 218
$wsComp.Range("D9").Value = 0

$wsComp.Range("A10").Value = "T4"
$wsComp.Range("B10").Value = "T5"
$wsComp.Range("C10").Value = 823
$wsComp.Range("D10").Value = 0

$wsComp.Range("E4").Select()

# --- Make comp_quantity_inst1 the active sheet/tab (activeTab=2) ---
$wsComp.Activate()
